$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H
$ws.Range("H1").Value = "FRA"

# FRA values for rows 2-23 (refactored out of another column into its own)
$values = @("Y","Y","Y","N","N","Y","Y","Y","U","N","N","U","Y","Y","Y","Y","Y","Y","Y","N","N","N")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}

# Update view: scroll position and selection to match final state
$ws.Range("H24").Select()
$excel.ActiveWindow.ScrollRow = 11
